$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$d = $ws.Cells.Item(2, 4)
$d.NumberFormat = '@'
$d.Value = '29.406.98'
$d.Style = 'Normal'
$ws.Cells.Item(2, 5).Value = '  -1.02%  '

# Row 3
$d = $ws.Cells.Item(3, 4)
$d.NumberFormat = '@'
$d.Value = '1.850.18'
$d.Style = 'Normal'
$ws.Cells.Item(3, 5).Value = '  -0.07%  '

# Row 4
$d = $ws.Cells.Item(4, 4)
$d.NumberFormat = '@'
$d.Value = '0.9992'
$d.Style = 'Normal'
$ws.Cells.Item(4, 5).Value = '  -0.07%  '

# Row 5
$d = $ws.Cells.Item(5, 4)
$d.NumberFormat = '@'
$d.Value = '242.09'
$d.Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -0.62%  '

# Row 6
$d = $ws.Cells.Item(6, 4)
$d.NumberFormat = '@'
$d.Value = '0.6293'
$d.Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -4.16%  '

# Row 7
$d = $ws.Cells.Item(7, 4)
$d.NumberFormat = '@'
$d.Value = '1.000'
$d.Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  -0.04%  '

# Row 8
$d = $ws.Cells.Item(8, 4)
$d.NumberFormat = '@'
$d.Value = '0.07615'
$d.Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  +1.60%  '

# Row 9
$d = $ws.Cells.Item(9, 4)
$d.NumberFormat = '@'
$d.Value = '0.2975'
$d.Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  -0.19%  '

# Row 10
$d = $ws.Cells.Item(10, 4)
$d.NumberFormat = '@'
$d.Value = '24.43'
$d.Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -0.43%  '

# Row 11
$ws.Cells.Item(11, 2).Value = 'WrappedEther'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$d = $ws.Cells.Item(11, 4)
$d.NumberFormat = '@'
$d.Value = '1.965.25'
$d.Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  +6.02%  '

# Row 12
$ws.Cells.Item(12, 2).Value = 'TRON'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$d = $ws.Cells.Item(12, 4)
$d.NumberFormat = '@'
$d.Value = '0.07721'
$d.Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  +0.96%  '

# Row 13
$d = $ws.Cells.Item(13, 4)
$d.NumberFormat = '@'
$d.Value = '5.007'
$d.Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  -0.78%  '

# Row 14
$d = $ws.Cells.Item(14, 4)
$d.NumberFormat = '@'
$d.Value = '0.6893'
$d.Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  +0.37%  '

# Row 15
$d = $ws.Cells.Item(15, 4)
$d.NumberFormat = '@'
$d.Value = '83.28'
$d.Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  -0.27%  '

# Row 16
$d = $ws.Cells.Item(16, 4)
$d.NumberFormat = '@'
$d.Value = '0.000009970'
$d.Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  +4.03%  '

# Row 17
$d = $ws.Cells.Item(17, 4)
$d.NumberFormat = '@'
$d.Value = '2.180.49'
$d.Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  +3.60%  '

# Row 18
$d = $ws.Cells.Item(18, 4)
$d.NumberFormat = '@'
$d.Value = '6.192'
$d.Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  +1.11%  '

# Row 19
$d = $ws.Cells.Item(19, 4)
$d.NumberFormat = '@'
$d.Value = '29.525.00'
$d.Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  -0.71%  '

# Row 20
$d = $ws.Cells.Item(20, 4)
$d.NumberFormat = '@'
$d.Value = '233.15'
$d.Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  -1.82%  '

# Row 21
$d = $ws.Cells.Item(21, 4)
$d.NumberFormat = '@'
$d.Value = '12.55'
$d.Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  -0.57%  '

# Row 22
$ws.Cells.Item(22, 2).Value = 'Chainlink'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$d = $ws.Cells.Item(22, 4)
$d.NumberFormat = '@'
$d.Value = '7.697'
$d.Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -0.39%  '

# Row 23
$ws.Cells.Item(23, 2).Value = 'Dai'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$d = $ws.Cells.Item(23, 4)
$d.NumberFormat = '@'
$d.Value = '1.000'
$d.Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +0.01%  '

# Row 24
$d = $ws.Cells.Item(24, 4)
$d.NumberFormat = '@'
$d.Value = '1.000'
$d.Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -0.09%  '

# Row 25
$d = $ws.Cells.Item(25, 4)
$d.NumberFormat = '@'
$d.Value = '155.03'
$d.Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -2.05%  '

# Row 26
$d = $ws.Cells.Item(26, 4)
$d.NumberFormat = '@'
$d.Value = '0.1395'
$d.Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  -2.17%  '

# Row 27
$d = $ws.Cells.Item(27, 4)
$d.NumberFormat = '@'
$d.Value = '8.488'
$d.Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -0.62%  '

# Row 28
$d = $ws.Cells.Item(28, 4)
$d.NumberFormat = '@'
$d.Value = '17.68'
$d.Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  -0.98%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  -1.22%  '

# Row 30
$d = $ws.Cells.Item(30, 4)
$d.NumberFormat = '@'
$d.Value = '0.05784'
$d.Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -4.30%  '

# Row 31
$d = $ws.Cells.Item(31, 4)
$d.NumberFormat = '@'
$d.Value = '1.256'
$d.Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  -1.32%  '

# Row 32
$d = $ws.Cells.Item(32, 4)
$d.NumberFormat = '@'
$d.Value = '4.131'
$d.Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -0.31%  '

# Row 33
$d = $ws.Cells.Item(33, 4)
$d.NumberFormat = '@'
$d.Value = '4.023'
$d.Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -1.16%  '

# Row 34
$d = $ws.Cells.Item(34, 4)
$d.NumberFormat = '@'
$d.Value = '1.883'
$d.Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +0.70%  '

# Row 35
$d = $ws.Cells.Item(35, 4)
$d.NumberFormat = '@'
$d.Value = '1.163'
$d.Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  -1.52%  '

# Row 36
$d = $ws.Cells.Item(36, 4)
$d.NumberFormat = '@'
$d.Value = '0.7209'
$d.Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  -0.80%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -0.55%  '

# Row 38
$d = $ws.Cells.Item(38, 4)
$d.NumberFormat = '@'
$d.Value = '1.247.15'
$d.Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  +3.64%  '

# Row 39
$d = $ws.Cells.Item(39, 4)
$d.NumberFormat = '@'
$d.Value = '2.794'
$d.Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -0.34%  '

# Row 40
$d = $ws.Cells.Item(40, 4)
$d.NumberFormat = '@'
$d.Value = '0.01805'
$d.Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  +0.83%  '

# Row 41
$d = $ws.Cells.Item(41, 4)
$d.NumberFormat = '@'
$d.Value = '0.9090'
$d.Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  -0.56%  '

# Row 42
$d = $ws.Cells.Item(42, 4)
$d.NumberFormat = '@'
$d.Value = '6.098'
$d.Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -2.99%  '

# Row 43
$d = $ws.Cells.Item(43, 4)
$d.NumberFormat = '@'
$d.Value = '2.087.23'
$d.Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +3.56%  '

# Row 44
$d = $ws.Cells.Item(44, 4)
$d.NumberFormat = '@'
$d.Value = '0.9993'
$d.Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  -0.06%  '

# Row 45
$d = $ws.Cells.Item(45, 4)
$d.NumberFormat = '@'
$d.Value = '67.90'
$d.Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  +1.93%  '

# Row 46
$d = $ws.Cells.Item(46, 4)
$d.NumberFormat = '@'
$d.Value = '101.78'
$d.Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  +0.71%  '

# Row 47
$d = $ws.Cells.Item(47, 4)
$d.NumberFormat = '@'
$d.Value = '7.293'
$d.Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -0.62%  '

# Row 48
$d = $ws.Cells.Item(48, 4)
$d.NumberFormat = '@'
$d.Value = '0.00000000121'
$d.Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  +1.87%  '

# Row 49
$d = $ws.Cells.Item(49, 4)
$d.NumberFormat = '@'
$d.Value = '9.199'
$d.Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  +0.78%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  -0.58%  '

# Row 51
$d = $ws.Cells.Item(51, 4)
$d.NumberFormat = '@'
$d.Value = '1.701'
$d.Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  +1.72%  '
